$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D:E").Insert()

Write-Host "dimension after insert done"
